$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2211.2
$ws.Range("I112").Value = 833
$ws.Range("J112").Value = 2399.1365
$ws.Range("K112").Value = 2499
$ws.Range("L112").Value = 7197.4095
$ws.Range("M112").Value = -1391
$ws.Range("N112").Value = -9413.4095

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1995.49
$ws.Range("I138").Value = 1252.0476
$ws.Range("J138").Value = 2193.114
$ws.Range("K138").Value = 3756.142800000001
$ws.Range("L138").Value = 6579.342000000001
$ws.Range("M138").Value = 1383.857199999999
$ws.Range("N138").Value = -16859.342

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15418.109
$ws.Range("I32").Value = 15857.266
$ws.Range("J32").Value = 11831.667
$ws.Range("K32").Value = 15857.266
$ws.Range("L32").Value = 11831.667
$ws.Range("M32").Value = -15570.266
$ws.Range("N32").Value = -12405.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1554.6666
$ws.Range("I61").Value = 1474
$ws.Range("J61").Value = 2200
$ws.Range("K61").Value = 1474
$ws.Range("L61").Value = 2200
$ws.Range("M61").Value = -1262
$ws.Range("N61").Value = -2624

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1641
$ws.Range("I110").Value = 1672.8572
$ws.Range("K110").Value = 1672.8572
$ws.Range("M110").Value = 372.1428000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1309.7084
$ws.Range("I132").Value = 1174.2727
$ws.Range("J132").Value = 2799.5
$ws.Range("K132").Value = 3522.8181
$ws.Range("L132").Value = 8398.5
$ws.Range("M132").Value = -992.8181
$ws.Range("N132").Value = -13458.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1554.6666
$ws.Range("I136").Value = 1474
$ws.Range("J136").Value = 2200
$ws.Range("K136").Value = 4422
$ws.Range("L136").Value = 6600
$ws.Range("M136").Value = -1872
$ws.Range("N136").Value = -11700

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 50000
$ws.Range("J40").Value = 50000
$ws.Range("L40").Value = 50000
$ws.Range("N40").Value = -50530

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 2314899.5
$ws.Range("I80").Value = 7407462
$ws.Range("J80").Value = 98.454544
$ws.Range("K80").Value = 7407462
$ws.Range("L80").Value = 98.454544
$ws.Range("M80").Value = -7406464
$ws.Range("N80").Value = -2094.454544

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 2314899.5
$ws.Range("I83").Value = 7407462
$ws.Range("J83").Value = 98.454544
$ws.Range("K83").Value = 37037310
$ws.Range("L83").Value = 492.27272
$ws.Range("M83").Value = -37032318
$ws.Range("N83").Value = -10476.27272

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2238
$ws.Range("I105").Value = 2238
$ws.Range("K105").Value = 2238
$ws.Range("M105").Value = -491

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2601.6
$ws.Range("I58").Value = 2185.3333
$ws.Range("J58").Value = 4266.6665
$ws.Range("K58").Value = 2185.3333
$ws.Range("L58").Value = 4266.6665
$ws.Range("M58").Value = -1982.3333
$ws.Range("N58").Value = -4672.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H93").Value = 19625
$ws.Range("I93").Value = 11000
$ws.Range("K93").Value = 11000
$ws.Range("M93").Value = -9128

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2601.6
$ws.Range("I136").Value = 2185.3333
$ws.Range("J136").Value = 4266.6665
$ws.Range("K136").Value = 6555.999899999999
$ws.Range("L136").Value = 12799.9995
$ws.Range("M136").Value = -4005.999899999999
$ws.Range("N136").Value = -17899.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 2375
$ws.Range("J54").Value = 2866.6667
$ws.Range("L54").Value = 8600.000100000001
$ws.Range("N54").Value = -9718.000100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1306.4762
$ws.Range("J68").Value = 1296
$ws.Range("L68").Value = 3888
$ws.Range("N68").Value = -5510

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1306.4762
$ws.Range("J71").Value = 1296
$ws.Range("L71").Value = 11664
$ws.Range("N71").Value = -19776

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 10500
$ws.Range("J80").Value = 4000
$ws.Range("L80").Value = 12000
$ws.Range("N80").Value = -13872

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 10500
$ws.Range("J83").Value = 4000
$ws.Range("L83").Value = 36000
$ws.Range("N83").Value = -45360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 886.2
$ws.Range("J131").Value = 912.9787
$ws.Range("L131").Value = 2738.9361
$ws.Range("N131").Value = -12818.9361

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1919
$ws.Range("I132").Value = 1099.2
$ws.Range("J132").Value = 2328.9
$ws.Range("K132").Value = 9892.800000000001
$ws.Range("L132").Value = 20960.1
$ws.Range("M132").Value = -7362.800000000001
$ws.Range("N132").Value = -26020.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 47622290
$ws.Range("J137").Value = 55559136
$ws.Range("L137").Value = 166677408
$ws.Range("N137").Value = -166687608

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 24545.455
$ws.Range("I4").Value = 20000
$ws.Range("J4").Value = 45000
$ws.Range("K4").Value = 20000
$ws.Range("L4").Value = 45000
$ws.Range("M4").Value = -19888
$ws.Range("N4").Value = -45224

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 30719.412
$ws.Range("I97").Value = 30719.412
$ws.Range("K97").Value = 30719.412
$ws.Range("M97").Value = -30223.412

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 12278.9
$ws.Range("J109").Value = 12278.9
$ws.Range("L109").Value = 12278.9
$ws.Range("N109").Value = -14358.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2250.2917
$ws.Range("I126").Value = 2077.5557
$ws.Range("K126").Value = 6232.6671
$ws.Range("M126").Value = -3762.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1880.4
$ws.Range("I93").Value = 1801
$ws.Range("J93").Value = 1999.5
$ws.Range("K93").Value = 1801
$ws.Range("L93").Value = 1999.5
$ws.Range("M93").Value = -553
$ws.Range("N93").Value = -4495.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4264
$ws.Range("I136").Value = 4718.3228
$ws.Range("K136").Value = 14154.9684
$ws.Range("M136").Value = -11604.9684

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 27528.5
$ws.Range("J64").Value = 27528.5
$ws.Range("L64").Value = 27528.5
$ws.Range("N64").Value = -28024.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H67").Value = 27528.5
$ws.Range("J67").Value = 27528.5
$ws.Range("L67").Value = 27528.5
$ws.Range("N67").Value = -29244.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 724.2857
$ws.Range("I96").Value = 775
$ws.Range("J96").Value = 656.6667
$ws.Range("K96").Value = 775
$ws.Range("L96").Value = 656.6667
$ws.Range("M96").Value = 598
$ws.Range("N96").Value = -3402.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 758.5517
$ws.Range("I113").Value = 413.61905
$ws.Range("J113").Value = 1664
$ws.Range("K113").Value = 1240.85715
$ws.Range("L113").Value = 4992
$ws.Range("M113").Value = 929.14285
$ws.Range("N113").Value = -9332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 16668316
$ws.Range("I122").Value = 35715576
$ws.Range("J122").Value = 1963.75
$ws.Range("K122").Value = 107146728
$ws.Range("L122").Value = 5891.25
$ws.Range("M122").Value = -107144278
$ws.Range("N122").Value = -10791.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1660.8518
$ws.Range("I132").Value = 1937.3334
$ws.Range("J132").Value = 693.1667
$ws.Range("K132").Value = 5812.0002
$ws.Range("L132").Value = 2079.5001
$ws.Range("M132").Value = -3282.0002
$ws.Range("N132").Value = -7139.5001
